# Update the Cxcl12-Ackr3 LR-pair sheet with the revised NATMI output
# (Dr Hou's advice): each of the 6 existing sender/target combinations is
# recomputed, and the 3 new sCs-sending rows are appended (rows 8-10), for
# a full 3x3 ECs/FAPs/sCs cross of sender x target clusters.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
  2 = @("ECs", "Cxcl12", "Ackr3", "ECs", 3, 1, "109.026058", "327.078174", "0.3049840938689738", "0.3049840938689738", 2, "0.6666666666666666", "11.497973", "34.493919", "0.09461878261088343", "0.09461878261088341", "1253.578671180434", "11282.20804062391", "0.0288572236775657", "0.02885722367756569")
  3 = @("ECs", "Cxcl12", "Ackr3", "FAPs", 3, 1, "109.026058", "327.078174", "0.3049840938689738", "0.3049840938689738", 3, 1, "99.43577566666666", "298.307327", "0.8182739724253053", "0.8182739724253052", "10841.09064510899", "97569.81580598089", "0.2495605460166974", "0.2495605460166974")
  4 = @("ECs", "Cxcl12", "Ackr3", "sCs", 3, 1, "109.026058", "327.078174", "0.3049840938689738", "0.3049840938689738", 3, 1, "10.585179", "31.755537", "0.0871072449638113", "0.0871072449638113", "1154.060339594382", "10386.54305634944", "0.02656632417471072", "0.02656632417471072")
  5 = @("FAPs", "Cxcl12", "Ackr3", "ECs", 3, 1, "89.97721833333333", "269.931655", "0.2516978134001918", "0.2516978134001917", 2, "0.6666666666666666", "11.497973", "34.493919", "0.09461878261088343", "0.09461878261088341", "1034.555627011772", "9311.000643105945", "0.02381534068974744", "0.02381534068974744")
  6 = @("FAPs", "Cxcl12", "Ackr3", "FAPs", 3, 1, "89.97721833333333", "269.931655", "0.2516978134001918", "0.2516978134001917", 3, 1, "99.43577566666666", "298.307327", "0.8182739724253053", "0.8182739724253052", "8946.954497304019", "80522.59047573617", "0.2059577696217381", "0.2059577696217381")
  7 = @("FAPs", "Cxcl12", "Ackr3", "sCs", 3, 1, "89.97721833333333", "269.931655", "0.2516978134001918", "0.2516978134001917", 3, 1, "10.585179", "31.755537", "0.0871072449638113", "0.0871072449638113", "952.4249619804148", "8571.824657823734", "0.02192470308870617", "0.02192470308870616")
  8 = @("sCs", "Cxcl12", "Ackr3", "ECs", 3, 1, "158.477852", "475.433556", "0.4433180927308344", "0.4433180927308344", 2, "0.6666666666666666", "11.497973", "34.493919", "0.09461878261088343", "0.09461878261088341", "1822.174063393996", "16399.56657054596", "0.04194621824357028", "0.04194621824357028")
  9 = @("sCs", "Cxcl12", "Ackr3", "FAPs", 3, 1, "158.477852", "475.433556", "0.4433180927308344", "0.4433180927308344", 3, 1, "99.43577566666666", "298.307327", "0.8182739724253053", "0.8182739724253052", "15758.3681396072", "141825.3132564648", "0.3627556567868698", "0.3627556567868697")
  10 = @("sCs", "Cxcl12", "Ackr3", "sCs", 3, 1, "158.477852", "475.433556", "0.4433180927308344", "0.4433180927308344", 3, 1, "10.585179", "31.755537", "0.0871072449638113", "0.0871072449638113", "1677.516430955508", "15097.64787859957", "0.03861621770039441", "0.0386162177003944")
}

foreach ($r in $rowData.Keys) {
  $vals = $rowData[$r]
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
  }
}
